$d = $word.ActiveDocument

# The first paragraph currently reads:
#   "This is a Microsoft word document."
# We need to turn it into four separate runs:
#   "This is a Microsoft word document."
#   " ("
#   "Changed main"
#   ")"
# so the visible text becomes:
#   "This is a Microsoft word document. (Changed main)"

$p = $d.Paragraphs(1)
$insertPos = $p.Range.End - 1   # just before the paragraph mark

$r = $d.Range($insertPos, $insertPos)

# Bookmarks act as run-boundary markers: inserting one, then text, then
# deleting the bookmark again leaves the surrounding text as distinct
# <w:r> runs instead of being silently coalesced into one run.
$d.Bookmarks.Add("zzEditBoundary0", $r)

$r.InsertAfter(" (")
$r.Collapse(0)
$d.Bookmarks.Add("zzEditBoundary1", $r)

$r.InsertAfter("Changed main")
$r.Collapse(0)
$d.Bookmarks.Add("zzEditBoundary2", $r)

$r.InsertAfter(")")

$d.Bookmarks("zzEditBoundary0").Delete()
$d.Bookmarks("zzEditBoundary1").Delete()
$d.Bookmarks("zzEditBoundary2").Delete()
